$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('Devin Booker', 'PG,SG', 'Phoenix Suns')
    ,@('Jamal Murray', 'PG,SG', 'Denver Nuggets')
    ,@('Trae Young', 'PG', 'Atlanta Hawks')
    ,@('LeBron James', 'SF,PF', 'Los Angeles Lakers')
    ,@('P.J. Washington', 'PF', 'Dallas Mavericks')
    ,@('Walker Kessler', 'C', 'Utah Jazz')
    ,@('Jalen Brunson', 'PG', 'New York Knicks')
    ,@('Immanuel Quickley', 'PG,SG', 'Toronto Raptors')
    ,@('Coby White', 'PG,SG', 'Chicago Bulls')
    ,@('Norman Powell', 'SG,SF', 'LA Clippers')
    ,@('Kawhi Leonard', 'SG,SF,PF', 'LA Clippers')
    ,@('Devin Vassell', 'SG,SF', 'San Antonio Spurs')
    ,@('Myles Turner', 'C', 'Indiana Pacers')
    ,@('Cole Anthony', 'PG', 'Orlando Magic')
    ,@('D''Angelo Russell', 'PG', 'Brooklyn Nets')
    ,@('Brandon Ingram', 'SG,SF,PF', 'New Orleans Pelicans')
    ,@('Jabari Smith Jr.', 'PF,C', 'Houston Rockets')
    ,@('Desmond Bane', 'SG,SF', 'Memphis Grizzlies')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
